$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.168.36'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.33%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.465.85'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.63%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '518.49'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.75%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.32'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.77%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('E8').Value = '  -1.86%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.474.05'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.44%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0980'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.38%  '

$ws.Range('E11').Value = '  -1.02%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.28'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.24%  '

$ws.Range('E13').Value = '  -2.68%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.905.23'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.62%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '58.066.60'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.43%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.89'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.84%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000134'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.46%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.474.26'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.45%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.45%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '318.69'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.60%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.15'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.19%  '

$ws.Range('E22').Value = '  -0.07%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.72'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.41%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.61'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.57%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.408'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.60%  '

$ws.Range('E26').Value = '  -0.20%  '

$ws.Range('E27').Value = '  -1.63%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.32'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.57%  '

$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '169.69'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.51%  '

$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0743'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.50%  '

$ws.Range('E31').Value = '  -1.92%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.68'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.39%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.17'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.80%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.02%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.04%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.98%  '

$ws.Range('E37').Value = '  -4.29%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.98'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.55%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.55'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.89%  '

$ws.Range('E40').Value = '  -3.48%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.792'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.72%  '

$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '272.86'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.94%  '

$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.42'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.35%  '

$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.27%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.592'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.73%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '122.97'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.15%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0905'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.01%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0487'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.52%  '

$ws.Range('E49').Value = '  -2.52%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.93'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.03%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.731.42'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.11%  '
